$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.339.91'
$ws.Range('E2').Value = '  +3.52%  '
$ws.Range('D3').Value = '3.640.56'
$ws.Range('E3').Value = '  +2.74%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '196.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +9.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '574.91'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('D7').Value = '3.635.22'
$ws.Range('E7').Value = '  +2.76%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.620'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.85%  '
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.679'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.156'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.79%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000297'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +19.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.11'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.42%  '
$ws.Range('D15').Value = '4.211.31'
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = '3.635.06'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('E17').Value = '  +0.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.71%  '
$ws.Range('D19').Value = '68.167.70'
$ws.Range('E19').Value = '  +3.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.81%  '
$ws.Range('E21').Value = '  +4.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '404.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.20'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +31.01%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.10'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.52%  '
$ws.Range('E26').Value = '  +4.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +5.12%  '
$ws.Range('E28').Value = '  +8.52%  '
$ws.Range('E29').Value = '  +1.76%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.25'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +24.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.23'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '32.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '689.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +16.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '12.28'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.48%  '
$ws.Range('E35').Value = '  +5.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '42.82'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.431'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +17.71%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').Value = '0.0₃0796'
$ws.Range('E40').Value = '  +10.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.95'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +24.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.141'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.29%  '
$ws.Range('E43').Value = '  +15.96%  '
$ws.Range('D44').Value = '3.221.68'
$ws.Range('E44').Value = '  +15.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.06'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +42.69%  '
$ws.Range('E46').Value = '  -0.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0422'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.54%  '
$ws.Range('E48').Value = '  +10.03%  '
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('E50').Value = '  +3.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '142.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.58%  '
